$wb = $excel.ActiveWorkbook

# Rename "Egresos" sheet to "Costos"
$wsCostos = $wb.Worksheets.Item("Egresos")
$wsCostos.Name = "Costos"

# Update selection on the "Inversiones" sheet
$wsInversiones = $wb.Worksheets.Item("Inversiones")
$wsInversiones.Range("I8:I13").Select()

# Activate the "Costos" sheet so it becomes the selected/active tab
$wsCostos.Activate()
$wsCostos.Range("D11").Select()
